# Driver sheet rewrite: align Scenario rows with the TestCaseName-driven
# layout used on the other pages, clear the ExecutionStatus column, and
# apply the same bordered/no-fill data style plus a bold+highlighted
# header style used elsewhere in the workbook.

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver")
$fmtSource = $wb.Worksheets.Item("ElementsPage")

# --- New data values for rows 2-5 (ExecutionStatus column cleared) -----
$ws.Range("A2").Value = "TC01"
$ws.Range("B2").Value = "Y"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "Create_form"
$ws.Range("E2").Value = "Create new Enrollment Form"

$ws.Range("A3").Value = "TC02"
$ws.Range("B3").Value = "N"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "Edit_form"
$ws.Range("E3").Value = "Edit new Enrollment Form"

$ws.Range("A4").Value = "TC03"
$ws.Range("B4").Value = "N"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = "Update_form"
$ws.Range("E4").Value = "Update new Enrollment Form"

$ws.Range("A5").Value = "TC04"
$ws.Range("B5").Value = "Y"
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = "Addnewrow_form"
$ws.Range("E5").Value = "Add new Row Enrollment Form"

# --- Re-style the data block (A2:E5) to match the bordered, unfilled ---
# --- style already used for data rows on the other sheets --------------
$fmtSource.Range("A2").Copy()
$ws.Range("A2:E5").PasteSpecial($xlPasteFormats)

# --- Re-style the header row (A1:E1): bold font + yellow fill + thin ---
# --- border. Start from the existing bold+bordered header style (blue --
# --- fill) used on the other sheets, then swap the fill to yellow so --
# --- the border/font definitions are reused verbatim and only the ----
# --- fill color needs a fresh style slot. -----------------------------
$fmtSource.Range("A1").Copy()
$ws.Range("A1:E1").PasteSpecial($xlPasteFormats)
$ws.Range("A1:E1").Interior.Color = 65535

$excel.CutCopyMode = $false

# --- Cosmetic bits also present in the target revision -----------------
$ws.PageSetup.Orientation = 1
$ws.Range("C2").Select()
